# Add a new development-diary row to the end of the table, mirroring the
# formatting of the preceding row (Rows.Add clones the last row's cell
# shading / widths / paragraph formatting automatically).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()

# Column 1: Date of Session
$cell1 = $newRow.Cells.Item(1)
$cell1.Range.Text = "21/12/2021"

# Column 2: Time Spent
$cell2 = $newRow.Cells.Item(2)
$cell2.Range.Text = "2 Hours 40 Minutes"

# Column 3: Development Segment
$cell3 = $newRow.Cells.Item(3)
$cell3.Range.Text = "World Generation – Objective 1"

# Column 4: Notes - multiple paragraphs (incl. a blank separator paragraph),
# inserted as raw OOXML so the exact run/paragraph structure is preserved.
$cell4 = $newRow.Cells.Item(4)
$cell4Xml = '<pkg:xmlData xmlns:pkg="http://schemas.openxmlformats.org/package/2006/content-types"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">Reworked various systems to improve memory-usage. First and foremost, this included removing the storing of integer values for each tile – instead deferring to using the enumerables that are generated from the integers. By removing the unnecessary storage </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>of the integer property values, each pixel in the map now takes up 16 bytes less. Additionally, the deciles system has been reworked to use the upper bounds and lower bounds as previously discussed, removing the need for the declaration of a new list fo</w:t></w:r><w:r><w:t>r calculating deciles.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>As a result of these changes, the memory usage of the features currently implemented has dropped by a full third. While the system overall still uses a heavy amount of memory, lowering this load by such a significant amount helps immensely.</w:t></w:r></w:p></pkg:xmlData>'
$cell4.Range.InsertXML($cell4Xml)
$placeholder = $cell4.Range.Paragraphs.Item(1)
$placeholder.Range.Delete()
